$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the existing header row (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new columns I (I0) and J (IF), rows 2-11
$values = @(
    @(1, 4),
    @(1, 4),
    @(2, 5),
    @(3, 6),
    @(9, 9),
    @(7, 8),
    @(2, 5),
    @(1, 4),
    @(1, 3),
    @(5, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
